$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-03 21:09:20"

$wsZhCn.Range("H2").Value = "2016-09-03 21:09:15"
$wsZhCn.Range("K2").Value = "2016-09-03 21:09:32"

$wsDeDe.Range("H2").Value = "2016-09-03 21:09:40"
